$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: I1 = "I0", J1 = "IF"
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the existing header formatting (bold font, border, centered) from H1
# onto the two new header cells so they reuse the same cell style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data columns I (I0) and J (IF) for rows 2-12.
$values = @(
    @(7, 7),
    @(7, 8),
    @(5, 7),
    @(6, 6),
    @(3, 4),
    @(8, 8),
    @(5, 5),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
